# Fix "Esther" Bible verse references in column A that have an
# erroneous trailing "16" suffix (e.g. "Esther 1:116" -> "Esther 1:1").
# A couple of rows (already correct, e.g. "Esther 2:13", "Esther 6:9")
# are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $value = $cell.Value2

    if ($value -ne $null -and $value -match '^(Esther \d+:\d+)16$') {
        $cell.Value2 = $matches[1]
    }
}
